$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 18518748
$ws.Cells.Item(6, 9).Value = 37037148
$ws.Cells.Item(6, 11).Value = 111111444
$ws.Cells.Item(6, 13).Value = -111111332
$ws.Cells.Item(8, 8).Value = 55555790
$ws.Cells.Item(8, 9).Value = 55555790
$ws.Cells.Item(8, 11).Value = 166667370
$ws.Cells.Item(8, 13).Value = -166667231
$ws.Cells.Item(9, 8).Value = 55555756
$ws.Cells.Item(9, 10).Value = 55555756
$ws.Cells.Item(9, 12).Value = 55555756
$ws.Cells.Item(9, 14).Value = -55556094
$ws.Cells.Item(61, 8).Value = 767.3333
$ws.Cells.Item(61, 9).Value = 767.3333
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 2301.9999
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -2129.9999
$ws.Cells.Item(61, 14).ClearContents()
$ws.Cells.Item(98, 8).Value = 67922
$ws.Cells.Item(98, 9).Value = 70162.25
$ws.Cells.Item(98, 11).Value = 70162.25
$ws.Cells.Item(98, 13).Value = -68664.25
$ws.Cells.Item(122, 8).Value = 67922
$ws.Cells.Item(122, 9).Value = 70162.25
$ws.Cells.Item(122, 11).Value = 210486.75
$ws.Cells.Item(122, 13).Value = -208036.75
$ws.Cells.Item(125, 8).Value = 3010.2
$ws.Cells.Item(125, 9).Value = 826
$ws.Cells.Item(125, 10).Value = 4466.3335
$ws.Cells.Item(125, 11).Value = 7434
$ws.Cells.Item(125, 12).Value = 40197.0015
$ws.Cells.Item(125, 13).Value = -4974
$ws.Cells.Item(125, 14).Value = -45117.0015
$ws.Cells.Item(138, 8).Value = 9259.467000000001
$ws.Cells.Item(138, 9).Value = 8416.666999999999
$ws.Cells.Item(138, 10).Value = 9821.333000000001
$ws.Cells.Item(138, 11).Value = 25250.001
$ws.Cells.Item(138, 12).Value = 29463.999
$ws.Cells.Item(138, 13).Value = -20110.001
$ws.Cells.Item(138, 14).Value = -39743.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 12991.667
$ws.Cells.Item(45, 9).Value = 17266.2
$ws.Cells.Item(45, 11).Value = 17266.2
$ws.Cells.Item(45, 13).Value = -16889.2
$ws.Cells.Item(61, 8).Value = 8613.714
$ws.Cells.Item(61, 9).Value = 10228.429
$ws.Cells.Item(61, 10).Value = 5384.2856
$ws.Cells.Item(61, 11).Value = 10228.429
$ws.Cells.Item(61, 12).Value = 5384.2856
$ws.Cells.Item(61, 13).Value = -10016.429
$ws.Cells.Item(61, 14).Value = -5808.2856
$ws.Cells.Item(136, 8).Value = 8613.714
$ws.Cells.Item(136, 9).Value = 10228.429
$ws.Cells.Item(136, 10).Value = 5384.2856
$ws.Cells.Item(136, 11).Value = 30685.287
$ws.Cells.Item(136, 12).Value = 16152.8568
$ws.Cells.Item(136, 13).Value = -28135.287
$ws.Cells.Item(136, 14).Value = -21252.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(70, 8).Value = 209555
$ws.Cells.Item(70, 10).Value = 209555
$ws.Cells.Item(70, 12).Value = 209555
$ws.Cells.Item(70, 14).Value = -210141
$ws.Cells.Item(73, 8).Value = 209555
$ws.Cells.Item(73, 10).Value = 209555
$ws.Cells.Item(73, 12).Value = 209555
$ws.Cells.Item(73, 14).Value = -211583

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2564.8298
$ws.Cells.Item(31, 9).Value = 1943.6111
$ws.Cells.Item(31, 10).Value = 2950.4138
$ws.Cells.Item(31, 11).Value = 1943.6111
$ws.Cells.Item(31, 12).Value = 2950.4138
$ws.Cells.Item(31, 13).Value = -1648.6111
$ws.Cells.Item(31, 14).Value = -3540.4138
$ws.Cells.Item(34, 8).Value = 2564.8298
$ws.Cells.Item(34, 9).Value = 1943.6111
$ws.Cells.Item(34, 10).Value = 2950.4138
$ws.Cells.Item(34, 11).Value = 1943.6111
$ws.Cells.Item(34, 12).Value = 2950.4138
$ws.Cells.Item(34, 13).Value = -1741.6111
$ws.Cells.Item(34, 14).Value = -3354.4138
$ws.Cells.Item(99, 8).Value = 2502000
$ws.Cells.Item(99, 9).Value = 2502000
$ws.Cells.Item(99, 11).Value = 2502000
$ws.Cells.Item(99, 13).Value = -2500502
$ws.Cells.Item(107, 8).Value = 62514480
$ws.Cells.Item(107, 9).Value = 83352010
$ws.Cells.Item(107, 11).Value = 83352010
$ws.Cells.Item(107, 13).Value = -83350090
$ws.Cells.Item(122, 8).Value = 1928.5714
$ws.Cells.Item(122, 9).Value = 1633.3334
$ws.Cells.Item(122, 11).Value = 4900.0002
$ws.Cells.Item(122, 13).Value = -2450.0002
$ws.Cells.Item(126, 8).Value = 2502000
$ws.Cells.Item(126, 9).Value = 2502000
$ws.Cells.Item(126, 11).Value = 7506000
$ws.Cells.Item(126, 13).Value = -7503530
$ws.Cells.Item(134, 8).Value = 1362748.5
$ws.Cells.Item(134, 9).Value = 1528486
$ws.Cells.Item(134, 10).Value = 3701.4
$ws.Cells.Item(134, 11).Value = 4585458
$ws.Cells.Item(134, 12).Value = 11104.2
$ws.Cells.Item(134, 13).Value = -4582923
$ws.Cells.Item(134, 14).Value = -16174.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 255.46666
$ws.Cells.Item(2, 9).Value = 212.16667
$ws.Cells.Item(2, 10).Value = 284.33334
$ws.Cells.Item(2, 11).Value = 1273.00002
$ws.Cells.Item(2, 12).Value = 1706.00004
$ws.Cells.Item(2, 13).Value = -1160.00002
$ws.Cells.Item(2, 14).Value = -1932.00004
$ws.Cells.Item(38, 8).Value = 1115.5667
$ws.Cells.Item(38, 9).Value = 170.81818
$ws.Cells.Item(38, 10).Value = 1662.5264
$ws.Cells.Item(38, 11).Value = 512.4545400000001
$ws.Cells.Item(38, 12).Value = 4987.5792
$ws.Cells.Item(38, 13).Value = -165.4545400000001
$ws.Cells.Item(38, 14).Value = -5681.5792
$ws.Cells.Item(60, 8).Value = 2042.5
$ws.Cells.Item(60, 9).Value = 463.75
$ws.Cells.Item(60, 11).Value = 1391.25
$ws.Cells.Item(60, 13).Value = -1140.25
$ws.Cells.Item(92, 8).Value = 598.1818
$ws.Cells.Item(92, 9).Value = 539.5
$ws.Cells.Item(92, 11).Value = 1618.5
$ws.Cells.Item(92, 13).Value = -370.5
$ws.Cells.Item(107, 8).Value = 1605.8334
$ws.Cells.Item(107, 9).Value = 716.6667
$ws.Cells.Item(107, 10).Value = 1783.6666
$ws.Cells.Item(107, 11).Value = 2150.0001
$ws.Cells.Item(107, 12).Value = 5350.9998
$ws.Cells.Item(107, 13).Value = -230.0001000000002
$ws.Cells.Item(107, 14).Value = -9190.9998
$ws.Cells.Item(125, 8).Value = 9460
$ws.Cells.Item(125, 9).Value = 2250
$ws.Cells.Item(125, 10).Value = 14266.667
$ws.Cells.Item(125, 11).Value = 6750
$ws.Cells.Item(125, 12).Value = 42800.001
$ws.Cells.Item(125, 13).Value = -1830
$ws.Cells.Item(125, 14).Value = -52640.001
$ws.Cells.Item(126, 8).Value = 27544
$ws.Cells.Item(126, 9).Value = 19500
$ws.Cells.Item(126, 11).Value = 58500
$ws.Cells.Item(126, 13).Value = -53560

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(38, 8).Value = 10000
$ws.Cells.Item(38, 10).Value = 10000
$ws.Cells.Item(38, 12).Value = 10000
$ws.Cells.Item(38, 14).Value = -10926
$ws.Cells.Item(97, 8).Value = 9195.959999999999
$ws.Cells.Item(97, 9).Value = 10727.75
$ws.Cells.Item(97, 10).Value = 3068.8
$ws.Cells.Item(97, 11).Value = 10727.75
$ws.Cells.Item(97, 12).Value = 3068.8
$ws.Cells.Item(97, 13).Value = -10231.75
$ws.Cells.Item(97, 14).Value = -4060.8
$ws.Cells.Item(102, 8).Value = 11730.333
$ws.Cells.Item(102, 9).Value = 12461.071
$ws.Cells.Item(102, 10).Value = 1500
$ws.Cells.Item(102, 11).Value = 12461.071
$ws.Cells.Item(102, 12).Value = 1500
$ws.Cells.Item(102, 13).Value = -10839.071
$ws.Cells.Item(102, 14).Value = -4744
$ws.Cells.Item(104, 8).Value = 50671
$ws.Cells.Item(104, 10).Value = 50671
$ws.Cells.Item(104, 12).Value = 50671
$ws.Cells.Item(104, 14).Value = -57659
$ws.Cells.Item(105, 8).Value = 81650
$ws.Cells.Item(105, 10).Value = 81650
$ws.Cells.Item(105, 12).Value = 81650
$ws.Cells.Item(105, 14).Value = -88638
$ws.Cells.Item(122, 8).Value = 24977.334
$ws.Cells.Item(122, 9).Value = 27098.5
$ws.Cells.Item(122, 10).Value = 8008
$ws.Cells.Item(122, 11).Value = 81295.5
$ws.Cells.Item(122, 12).Value = 24024
$ws.Cells.Item(122, 13).Value = -78845.5
$ws.Cells.Item(122, 14).Value = -28924

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 51719.555
$ws.Cells.Item(7, 9).Value = 63711
$ws.Cells.Item(7, 10).Value = 9749.5
$ws.Cells.Item(7, 11).Value = 63711
$ws.Cells.Item(7, 12).Value = 9749.5
$ws.Cells.Item(7, 13).Value = -63599
$ws.Cells.Item(7, 14).Value = -9973.5
$ws.Cells.Item(40, 8).Value = 66957.69500000001
$ws.Cells.Item(40, 9).Value = 84551.78
$ws.Cells.Item(40, 10).Value = 27371
$ws.Cells.Item(40, 11).Value = 84551.78
$ws.Cells.Item(40, 12).Value = 27371
$ws.Cells.Item(40, 13).Value = -84415.78
$ws.Cells.Item(40, 14).Value = -27643
$ws.Cells.Item(100, 8).Value = 5526.0527
$ws.Cells.Item(100, 9).Value = 5764.4116
$ws.Cells.Item(100, 11).Value = 5764.4116
$ws.Cells.Item(100, 13).Value = -5223.4116
$ws.Cells.Item(105, 8).Value = 110615
$ws.Cells.Item(105, 10).Value = 110615
$ws.Cells.Item(105, 12).Value = 110615
$ws.Cells.Item(105, 14).Value = -117603
$ws.Cells.Item(122, 8).Value = 5230.6924
$ws.Cells.Item(122, 9).Value = 3399.75
$ws.Cells.Item(122, 11).Value = 10199.25
$ws.Cells.Item(122, 13).Value = -7749.25
$ws.Cells.Item(126, 8).Value = 51719.555
$ws.Cells.Item(126, 9).Value = 63711
$ws.Cells.Item(126, 10).Value = 9749.5
$ws.Cells.Item(126, 11).Value = 191133
$ws.Cells.Item(126, 12).Value = 29248.5
$ws.Cells.Item(126, 13).Value = -188663
$ws.Cells.Item(126, 14).Value = -34188.5
$ws.Cells.Item(136, 8).Value = 3881.5
$ws.Cells.Item(136, 9).Value = 1971.3
$ws.Cells.Item(136, 10).Value = 4942.722
$ws.Cells.Item(136, 11).Value = 5913.9
$ws.Cells.Item(136, 12).Value = 14828.166
$ws.Cells.Item(136, 13).Value = -3363.9
$ws.Cells.Item(136, 14).Value = -19928.166

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 19423.363
$ws.Cells.Item(81, 9).Value = 20715.7
$ws.Cells.Item(81, 11).Value = 41431.4
$ws.Cells.Item(81, 13).Value = -40370.4
$ws.Cells.Item(84, 8).Value = 19423.363
$ws.Cells.Item(84, 9).Value = 20715.7
$ws.Cells.Item(84, 11).Value = 207157
$ws.Cells.Item(84, 13).Value = -201853
$ws.Cells.Item(122, 8).Value = 7762.05
$ws.Cells.Item(122, 9).Value = 5520.5
$ws.Cells.Item(122, 11).Value = 16561.5
$ws.Cells.Item(122, 13).Value = -14111.5
$ws.Cells.Item(136, 8).Value = 3887.7632
$ws.Cells.Item(136, 9).Value = 2691.037
$ws.Cells.Item(136, 11).Value = 8073.110999999999
$ws.Cells.Item(136, 13).Value = -5523.110999999999
